$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J65").Value = "2023-12-14T12:59:00"
